$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab name) to reflect the new "through" date
$ws.Name = "Through 2021-12-12"

# Update the December label cell (A13) to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-12)"

# Update H12 (2021 November value)
$ws.Range("H12").Value = 200

# Update row 13 (December values for years 2015-2021, columns B-H)
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 34
$ws.Range("D13").Value = 40
$ws.Range("E13").Value = 27
$ws.Range("F13").Value = 18
$ws.Range("G13").Value = 60
$ws.Range("H13").Value = 90

# Update row 14 (Total values for years 2015-2021, columns B-H)
$ws.Range("B14").Value = 303
$ws.Range("C14").Value = 597
$ws.Range("D14").Value = 861
$ws.Range("E14").Value = 709
$ws.Range("F14").Value = 552
$ws.Range("G14").Value = 1324
$ws.Range("H14").Value = 1734
